$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 481.75
$ws.Range("I8").Value = 845.25
$ws.Range("K8").Value = 2535.75
$ws.Range("M8").Value = -2396.75
$ws.Range("H98").Value = 2447.3264
$ws.Range("I98").Value = 2591.7441
$ws.Range("J98").Value = 1412.3334
$ws.Range("K98").Value = 2591.7441
$ws.Range("L98").Value = 1412.3334
$ws.Range("M98").Value = -1093.7441
$ws.Range("N98").Value = -4408.3334
$ws.Range("H100").Value = 1948.1111
$ws.Range("I100").Value = 1944.375
$ws.Range("K100").Value = 1944.375
$ws.Range("M100").Value = -1403.375
$ws.Range("H116").Value = 2899.8
$ws.Range("I116").Value = 2749.5
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2749.5
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 692.5
$ws.Range("N116").Value = -9884
$ws.Range("H122").Value = 2447.3264
$ws.Range("I122").Value = 2591.7441
$ws.Range("J122").Value = 1412.3334
$ws.Range("K122").Value = 7775.2323
$ws.Range("L122").Value = 4237.0002
$ws.Range("M122").Value = -5325.2323
$ws.Range("N122").Value = -9137.0002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 852.5625
$ws.Range("I2").Value = 688.5
$ws.Range("J2").Value = 1016.625
$ws.Range("K2").Value = 688.5
$ws.Range("L2").Value = 1016.625
$ws.Range("M2").Value = -575.5
$ws.Range("N2").Value = -1242.625
$ws.Range("H32").Value = 6655.1777
$ws.Range("I32").Value = 5439.225
$ws.Range("K32").Value = 5439.225
$ws.Range("M32").Value = -5152.225
$ws.Range("H45").Value = 1004.4839
$ws.Range("I45").Value = 923.6667
$ws.Range("K45").Value = 923.6667
$ws.Range("M45").Value = -546.6667
$ws.Range("H74").Value = 1909.421
$ws.Range("I74").Value = 1633.6875
$ws.Range("K74").Value = 1633.6875
$ws.Range("M74").Value = -759.6875
$ws.Range("H77").Value = 1909.421
$ws.Range("I77").Value = 1633.6875
$ws.Range("K77").Value = 8168.4375
$ws.Range("M77").Value = -3800.4375
$ws.Range("H110").Value = 1327.0454
$ws.Range("I110").Value = 883.5
$ws.Range("J110").Value = 3323
$ws.Range("K110").Value = 883.5
$ws.Range("L110").Value = 3323
$ws.Range("M110").Value = 1161.5
$ws.Range("N110").Value = -7413
$ws.Range("H116").Value = 852.5625
$ws.Range("I116").Value = 688.5
$ws.Range("J116").Value = 1016.625
$ws.Range("K116").Value = 688.5
$ws.Range("L116").Value = 1016.625
$ws.Range("M116").Value = 1605.5
$ws.Range("N116").Value = -5604.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 852.5625
$ws.Range("I3").Value = 688.5
$ws.Range("J3").Value = 1016.625
$ws.Range("K3").Value = 688.5
$ws.Range("L3").Value = 1016.625
$ws.Range("M3").Value = -574.5
$ws.Range("N3").Value = -1244.625
$ws.Range("H25").Value = 666
$ws.Range("I25").Value = 655
$ws.Range("J25").Value = 688
$ws.Range("K25").Value = 655
$ws.Range("L25").Value = 688
$ws.Range("M25").Value = -420
$ws.Range("N25").Value = -1158
$ws.Range("H80").Value = 480.05884
$ws.Range("I80").Value = 66
$ws.Range("J80").Value = 607.46155
$ws.Range("K80").Value = 66
$ws.Range("L80").Value = 607.46155
$ws.Range("M80").Value = 932
$ws.Range("N80").Value = -2603.46155
$ws.Range("H83").Value = 480.05884
$ws.Range("I83").Value = 66
$ws.Range("J83").Value = 607.46155
$ws.Range("K83").Value = 330
$ws.Range("L83").Value = 3037.30775
$ws.Range("M83").Value = 4662
$ws.Range("N83").Value = -13021.30775
$ws.Range("H134").Value = 5127
$ws.Range("I134").Value = 970.9167
$ws.Range("K134").Value = 2912.7501
$ws.Range("M134").Value = -377.7501000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 58824900
$ws.Range("I16").Value = 90910390
$ws.Range("J16").Value = 1496.6666
$ws.Range("K16").Value = 90910390
$ws.Range("L16").Value = 1496.6666
$ws.Range("M16").Value = -90910103
$ws.Range("N16").Value = -2070.6666
$ws.Range("H31").Value = 1283.88
$ws.Range("I31").Value = 1221.8823
$ws.Range("K31").Value = 1221.8823
$ws.Range("M31").Value = -926.8823
$ws.Range("H34").Value = 1283.88
$ws.Range("I34").Value = 1221.8823
$ws.Range("K34").Value = 1221.8823
$ws.Range("M34").Value = -1019.8823
$ws.Range("H107").Value = 771
$ws.Range("I107").Value = 383.26666
$ws.Range("K107").Value = 383.26666
$ws.Range("M107").Value = 1536.73334
$ws.Range("H113").Value = 58824900
$ws.Range("I113").Value = 90910390
$ws.Range("J113").Value = 1496.6666
$ws.Range("K113").Value = 90910390
$ws.Range("L113").Value = 1496.6666
$ws.Range("M113").Value = -90908220
$ws.Range("N113").Value = -5836.6666
$ws.Range("H132").Value = 1878.6666
$ws.Range("I132").Value = 1587.3158
$ws.Range("J132").Value = 2274.0715
$ws.Range("K132").Value = 4761.9474
$ws.Range("L132").Value = 6822.2145
$ws.Range("M132").Value = -2231.9474
$ws.Range("N132").Value = -11882.2145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 246340.69
$ws.Range("I4").Value = 99936.55499999999
$ws.Range("J4").Value = 575750
$ws.Range("K4").Value = 299809.665
$ws.Range("L4").Value = 1727250
$ws.Range("M4").Value = -299697.665
$ws.Range("N4").Value = -1727474
$ws.Range("H29").Value = 450.2857
$ws.Range("I29").Value = 83.333336
$ws.Range("K29").Value = 250.000008
$ws.Range("M29").Value = 26.99999199999999
$ws.Range("H86").Value = 354.22223
$ws.Range("I86").Value = 256.5
$ws.Range("J86").Value = 549.6667
$ws.Range("K86").Value = 769.5
$ws.Range("L86").Value = 1649.0001
$ws.Range("M86").Value = 416.5
$ws.Range("N86").Value = -4021.0001
$ws.Range("H89").Value = 354.22223
$ws.Range("I89").Value = 256.5
$ws.Range("J89").Value = 549.6667
$ws.Range("K89").Value = 2308.5
$ws.Range("L89").Value = 4947.0003
$ws.Range("M89").Value = 3619.5
$ws.Range("N89").Value = -16803.0003
$ws.Range("H131").Value = 10754938
$ws.Range("I131").Value = 111111384
$ws.Range("J131").Value = 2461.3572
$ws.Range("K131").Value = 333334152
$ws.Range("L131").Value = 7384.071599999999
$ws.Range("M131").Value = -333329112
$ws.Range("N131").Value = -17464.0716
$ws.Range("H137").Value = 25003504
$ws.Range("I137").Value = 125001624
$ws.Range("J137").Value = 3974.4167
$ws.Range("K137").Value = 375004872
$ws.Range("L137").Value = 11923.2501
$ws.Range("M137").Value = -374999772
$ws.Range("N137").Value = -22123.2501

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H102").Value = 1057.3572
$ws.Range("I102").Value = 724.0833
$ws.Range("J102").Value = 3057
$ws.Range("K102").Value = 724.0833
$ws.Range("L102").Value = 3057
$ws.Range("M102").Value = 897.9167
$ws.Range("N102").Value = -6301
